$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(20).Insert()
$ws.Range("A20").Value = 4
$ws.Range("B20").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C20").Value = "Los Lagos"
$ws.Range("D20").Value = 44771
$ws.Range("E20").Value = 10
$ws.Range("F20").Value = 100112012
$ws.Range("G20").Value = "Espinaca"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 40
$ws.Range("K20").Value = 14000
$ws.Range("L20").Value = 14000
$ws.Range("M20").Value = 14000
$ws.Range("N20").Value = "$/cuna 10 kilos"
$ws.Range("O20").Value = "Región Metropolitana"
$ws.Range("P20").Value = 1400
$ws.Range("Q20").Value = 10
$ws.Range("R20").Value = "Hortaliza"
